$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the original data (Name Center -> Description) before editing ---
# Original rows 1..9: column A = "xxx Center", column B = long description (definition text)
$origDescByFull = @{}
for ($r = 1; $r -le 9; $r++) {
    $full = $ws.Cells.Item($r, 1).Value2
    $desc = $ws.Cells.Item($r, 2).Value2
    $origDescByFull[$full] = $desc
}

# --- Clear the whole sheet so we can rebuild it cleanly ---
$ws.Cells.ClearContents()

# --- Step 1: write the description text into column B, rows 2-10, in original order ---
# (Head, Ajna, Throat, G, Ego, Solar Plexus, Sacral, Spleen, Root)
# This establishes shared-string slots 0-8 in that order.
$ws.Cells.Item(2, 2).Value = $origDescByFull["Head Center"]
$ws.Cells.Item(3, 2).Value = $origDescByFull["Ajna Center"]
$ws.Cells.Item(4, 2).Value = $origDescByFull["Throat Center"]
$ws.Cells.Item(5, 2).Value = $origDescByFull["G Center"]
$ws.Cells.Item(6, 2).Value = $origDescByFull["Ego Center"]
$ws.Cells.Item(7, 2).Value = $origDescByFull["Solar Plexus Center"]
$ws.Cells.Item(8, 2).Value = $origDescByFull["Sacral Center"]
$ws.Cells.Item(9, 2).Value = $origDescByFull["Spleen Center"]
$ws.Cells.Item(10, 2).Value = $origDescByFull["Root Center"]

# --- Step 2: write the header row (establishes shared-string slots 9-13) ---
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Description_def"
$ws.Cells.Item(1, 3).Value = "Description_undef"
$ws.Cells.Item(1, 4).Value = "Description_def_FR"
$ws.Cells.Item(1, 5).Value = "Description_undef_FR"

# --- Step 3: write the short center names into column A (establishes shared-string slots 14-22) ---
# Order: Head, Root, Ego, Spleen, Solar Plexus, Sacral, G, Throat, Ajna
$ws.Cells.Item(2, 1).Value = "Head"
$ws.Cells.Item(10, 1).Value = "Root"
$ws.Cells.Item(6, 1).Value = "Ego"
$ws.Cells.Item(9, 1).Value = "Spleen"
$ws.Cells.Item(7, 1).Value = "Solar Plexus"
$ws.Cells.Item(8, 1).Value = "Sacral"
$ws.Cells.Item(5, 1).Value = "G"
$ws.Cells.Item(4, 1).Value = "Throat"
$ws.Cells.Item(3, 1).Value = "Ajna"

# --- Formatting: bold header row ---
$ws.Range("A1:E1").Font.Bold = $true

# --- Column A width ---
$ws.Columns("A:A").ColumnWidth = 28.75

# --- Undo the automatic row-height expansion triggered by multi-line cell text ---
$ws.Rows("1:10").AutoFit()

# --- Selection matches the final state ---
$ws.Range("A4").Select()
